# Update natural products (np_harvest_individuals / NP_quotas / NP_Struck and
# Lost rate) metadata notes on the Sheet1 "Meta_Data" table to reflect the
# Pacific walrus PBR reference and the corrected Narwhal struck-and-lost
# multiplier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (np_harvest_individuals_arc2016) - note more detail on the years of
# harp/hooded seal + pacific walrus catch data available.
$ws.Range("C9").Value = "Harp and Hooded Seal catches. Contacted person about pacific walrus catch - have 5 years of data included 2006-2010"

# Row 10 (NP_quotas) - note that the Pacific Walrus PBR is a joint US/Russia
# limit, and that there are no hard quotas for seals/walrus (PBR is used
# instead as a sustainable-limit substitute).
$ws.Range("C10").Value = "ICES Quotas for Harp and Hooded Seals. Pacific Walrus PBR - joint with USA"
$ws.Range("G10").Value = "No quotas set for seals or walrus. PBR for pacific walrus used as sustainable limit from here https://www.fws.gov/alaska/fisheries/mmm/stock/Revised_April_2014_Pacific_Walrus_SAR.pdf"

# Row 11 (NP_Struck and Lost rate) - fix the narwhal multiplier (1.42 -> 1.28)
# and clarify why the walrus figure keeps 1.42 (PBR, not a quota).
$ws.Range("E11").Value = "Narwhal 28% struck and lost - multiply by 1.28"
$ws.Range("G11").Value = "Walrus 42% struck and lost - figures multiplied by 1.42 as reference is PBR not quota"

# Reflect the edits in the sheet's current view: scroll the frozen pane back
# up so row 9 is visible again, and leave the active selection on C10 (the
# NP_quotas note that was just revised).
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
$ws.Range("C10").Select()
